# Regenerate the "K" column (column G) of the save_data sheet with
# recalculated strike-count values (K) in place of the old Strike# values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value (column G), as recomputed by the regen script.
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 2
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 0
    36 = 0
    37 = 1
    38 = 1
    39 = 0
    40 = 4
    41 = 1
    42 = 1
    43 = 2
    44 = 0
    45 = 1
    46 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
